$p = $ppt.ActivePresentation
$s = $p.Slides.Add(2, 12)

# ---- Group 1 (rot=0) ----
$rect1 = $s.Shapes.AddShape(1, 0, 0, 1, 1)
$rect1.Left = 322.9908267716535
$rect1.Top = 276.0000393700788
$rect1.Width = 19.0092519685039
$rect1.Height = 15.2810629921260
$rect1.Name = "Rectangle 4"
$rect1.Fill.Visible = 0
$rect1.Line.ForeColor.ObjectThemeColor = 1
$rect1.Line.Weight = 1
$rect1.Line.DashStyle = 4
$rect1.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$oval1 = $s.Shapes.AddShape(9, 0, 0, 1, 1)
$oval1.Left = 252.0000393700788
$oval1.Top = 222.0000393700788
$oval1.Width = 141.9816929133858
$oval1.Height = 138.5621653543307
$oval1.Name = "Oval 5"
$oval1.Fill.Visible = 0
$oval1.Line.ForeColor.ObjectThemeColor = 1
$oval1.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$conn1_1 = $s.Shapes.AddConnector(1, 0, 0, 1, 1)
$conn1_1.Left = 322.9908267716535
$conn1_1.Top = 221.1108267716535
$conn1_1.Width = 96.0000393700787
$conn1_1.Height = 0.0000393700787
$conn1_1.Name = "Straight Arrow Connector 6"
$conn1_1.Line.Weight = 3
$conn1_1.Line.ForeColor.RGB = RGB(255,0,0)
$conn1_1.Line.EndArrowheadStyle = 3

$conn2_1 = $s.Shapes.AddConnector(1, 0, 0, 1, 1)
$conn2_1.Left = 394.9346850393701
$conn2_1.Top = 291.2810629921260
$conn2_1.Width = 0.0000393700787
$conn2_1.Height = 96.0000393700787
$conn2_1.Name = "Straight Arrow Connector 7"
$conn2_1.Line.Weight = 3
$conn2_1.Line.ForeColor.RGB = RGB(255,0,0)
$conn2_1.Line.EndArrowheadStyle = 3

$conn3_1 = $s.Shapes.AddLine(0, 0, 1, 1)
$conn3_1.Left = 322.9908267716535
$conn3_1.Top = 222.0000393700788
$conn3_1.Width = 0.0001181102362
$conn3_1.Height = 69.2810629921260
$conn3_1.Flip(0)
$conn3_1.Name = "Straight Connector 10"
$conn3_1.Line.ForeColor.ObjectThemeColor = 1
$conn3_1.Line.Weight = 1
$conn3_1.Line.DashStyle = 4
$conn3_1.ConnectorFormat.BeginConnect($oval1, 0)

$conn4_1 = $s.Shapes.AddLine(0, 0, 1, 1)
$conn4_1.Left = 329.9727165354331
$conn4_1.Top = 291.7107480314961
$conn4_1.Width = 67.9358661417323
$conn4_1.Height = 0.0000393700787
$conn4_1.Name = "Straight Connector 11"
$conn4_1.Line.ForeColor.ObjectThemeColor = 1
$conn4_1.Line.Weight = 1
$conn4_1.Line.DashStyle = 4

$grp1 = $s.Shapes.Range(@($rect1.Name, $oval1.Name, $conn1_1.Name, $conn2_1.Name, $conn3_1.Name, $conn4_1.Name)).Group()
$grp1.Name = "Group 3"
$grp1.Left = 114.0000393700787
$grp1.Top = 120.0000393700787
$grp1.Width = 120.0000393700787
$grp1.Height = 120.7883858267717

# ---- Group 2 (rot=90) ----
$rect2 = $s.Shapes.AddShape(1, 0, 0, 1, 1)
$rect2.Left = 322.9908267716535
$rect2.Top = 276.0000393700788
$rect2.Width = 19.0092519685039
$rect2.Height = 15.2810629921260
$rect2.Name = "Rectangle 27"
$rect2.Fill.Visible = 0
$rect2.Line.ForeColor.ObjectThemeColor = 1
$rect2.Line.Weight = 1
$rect2.Line.DashStyle = 4
$rect2.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$oval2 = $s.Shapes.AddShape(9, 0, 0, 1, 1)
$oval2.Left = 252.0000393700788
$oval2.Top = 222.0000393700788
$oval2.Width = 141.9816929133858
$oval2.Height = 138.5621653543307
$oval2.Name = "Oval 28"
$oval2.Fill.Visible = 0
$oval2.Line.ForeColor.ObjectThemeColor = 1
$oval2.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$conn1_2 = $s.Shapes.AddConnector(1, 0, 0, 1, 1)
$conn1_2.Left = 322.9908267716535
$conn1_2.Top = 221.1108267716535
$conn1_2.Width = 96.0000393700787
$conn1_2.Height = 0.0000393700787
$conn1_2.Name = "Straight Arrow Connector 29"
$conn1_2.Line.Weight = 3
$conn1_2.Line.ForeColor.RGB = RGB(255,0,0)
$conn1_2.Line.EndArrowheadStyle = 3

$conn2_2 = $s.Shapes.AddConnector(1, 0, 0, 1, 1)
$conn2_2.Left = 394.9346850393701
$conn2_2.Top = 291.2810629921260
$conn2_2.Width = 0.0000393700787
$conn2_2.Height = 96.0000393700787
$conn2_2.Name = "Straight Arrow Connector 30"
$conn2_2.Line.Weight = 3
$conn2_2.Line.ForeColor.RGB = RGB(255,0,0)
$conn2_2.Line.EndArrowheadStyle = 3

$conn3_2 = $s.Shapes.AddLine(0, 0, 1, 1)
$conn3_2.Left = 322.9908267716535
$conn3_2.Top = 222.0000393700788
$conn3_2.Width = 0.0001181102362
$conn3_2.Height = 69.2810629921260
$conn3_2.Flip(0)
$conn3_2.Name = "Straight Connector 31"
$conn3_2.Line.ForeColor.ObjectThemeColor = 1
$conn3_2.Line.Weight = 1
$conn3_2.Line.DashStyle = 4
$conn3_2.ConnectorFormat.BeginConnect($oval2, 0)

$conn4_2 = $s.Shapes.AddLine(0, 0, 1, 1)
$conn4_2.Left = 329.9727165354331
$conn4_2.Top = 291.7107480314961
$conn4_2.Width = 67.9358661417323
$conn4_2.Height = 0.0000393700787
$conn4_2.Name = "Straight Connector 32"
$conn4_2.Line.ForeColor.ObjectThemeColor = 1
$conn4_2.Line.Weight = 1
$conn4_2.Line.DashStyle = 4

$grp2 = $s.Shapes.Range(@($rect2.Name, $oval2.Name, $conn1_2.Name, $conn2_2.Name, $conn3_2.Name, $conn4_2.Name)).Group()
$grp2.Name = "Group 26"
$grp2.Left = 396.0000393700788
$grp2.Top = 120.0000393700787
$grp2.Width = 120.0000393700787
$grp2.Height = 120.7883858267717
$grp2.Rotation = 90

# ---- Group 3 (rot=180) ----
$rect3 = $s.Shapes.AddShape(1, 0, 0, 1, 1)
$rect3.Left = 322.9908267716535
$rect3.Top = 276.0000393700788
$rect3.Width = 19.0092519685039
$rect3.Height = 15.2810629921260
$rect3.Name = "Rectangle 34"
$rect3.Fill.Visible = 0
$rect3.Line.ForeColor.ObjectThemeColor = 1
$rect3.Line.Weight = 1
$rect3.Line.DashStyle = 4
$rect3.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$oval3 = $s.Shapes.AddShape(9, 0, 0, 1, 1)
$oval3.Left = 252.0000393700788
$oval3.Top = 222.0000393700788
$oval3.Width = 141.9816929133858
$oval3.Height = 138.5621653543307
$oval3.Name = "Oval 35"
$oval3.Fill.Visible = 0
$oval3.Line.ForeColor.ObjectThemeColor = 1
$oval3.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$conn1_3 = $s.Shapes.AddConnector(1, 0, 0, 1, 1)
$conn1_3.Left = 322.9908267716535
$conn1_3.Top = 221.1108267716535
$conn1_3.Width = 96.0000393700787
$conn1_3.Height = 0.0000393700787
$conn1_3.Name = "Straight Arrow Connector 36"
$conn1_3.Line.Weight = 3
$conn1_3.Line.ForeColor.RGB = RGB(255,0,0)
$conn1_3.Line.EndArrowheadStyle = 3

$conn2_3 = $s.Shapes.AddConnector(1, 0, 0, 1, 1)
$conn2_3.Left = 394.9346850393701
$conn2_3.Top = 291.2810629921260
$conn2_3.Width = 0.0000393700787
$conn2_3.Height = 96.0000393700787
$conn2_3.Name = "Straight Arrow Connector 37"
$conn2_3.Line.Weight = 3
$conn2_3.Line.ForeColor.RGB = RGB(255,0,0)
$conn2_3.Line.EndArrowheadStyle = 3

$conn3_3 = $s.Shapes.AddLine(0, 0, 1, 1)
$conn3_3.Left = 322.9908267716535
$conn3_3.Top = 222.0000393700788
$conn3_3.Width = 0.0001181102362
$conn3_3.Height = 69.2810629921260
$conn3_3.Flip(0)
$conn3_3.Name = "Straight Connector 38"
$conn3_3.Line.ForeColor.ObjectThemeColor = 1
$conn3_3.Line.Weight = 1
$conn3_3.Line.DashStyle = 4
$conn3_3.ConnectorFormat.BeginConnect($oval3, 0)

$conn4_3 = $s.Shapes.AddLine(0, 0, 1, 1)
$conn4_3.Left = 329.9727165354331
$conn4_3.Top = 291.7107480314961
$conn4_3.Width = 67.9358661417323
$conn4_3.Height = 0.0000393700787
$conn4_3.Name = "Straight Connector 39"
$conn4_3.Line.ForeColor.ObjectThemeColor = 1
$conn4_3.Line.Weight = 1
$conn4_3.Line.DashStyle = 4

$grp3 = $s.Shapes.Range(@($rect3.Name, $oval3.Name, $conn1_3.Name, $conn2_3.Name, $conn3_3.Name, $conn4_3.Name)).Group()
$grp3.Name = "Group 33"
$grp3.Left = 105.0142125984252
$grp3.Top = 336.0000393700788
$grp3.Width = 120.0000393700787
$grp3.Height = 120.7883858267717
$grp3.Rotation = 180

# ---- Group 4 (rot=270) ----
$rect4 = $s.Shapes.AddShape(1, 0, 0, 1, 1)
$rect4.Left = 322.9908267716535
$rect4.Top = 276.0000393700788
$rect4.Width = 19.0092519685039
$rect4.Height = 15.2810629921260
$rect4.Name = "Rectangle 41"
$rect4.Fill.Visible = 0
$rect4.Line.ForeColor.ObjectThemeColor = 1
$rect4.Line.Weight = 1
$rect4.Line.DashStyle = 4
$rect4.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$oval4 = $s.Shapes.AddShape(9, 0, 0, 1, 1)
$oval4.Left = 252.0000393700788
$oval4.Top = 222.0000393700788
$oval4.Width = 141.9816929133858
$oval4.Height = 138.5621653543307
$oval4.Name = "Oval 42"
$oval4.Fill.Visible = 0
$oval4.Line.ForeColor.ObjectThemeColor = 1
$oval4.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$conn1_4 = $s.Shapes.AddConnector(1, 0, 0, 1, 1)
$conn1_4.Left = 322.9908267716535
$conn1_4.Top = 221.1108267716535
$conn1_4.Width = 96.0000393700787
$conn1_4.Height = 0.0000393700787
$conn1_4.Name = "Straight Arrow Connector 43"
$conn1_4.Line.Weight = 3
$conn1_4.Line.ForeColor.RGB = RGB(255,0,0)
$conn1_4.Line.EndArrowheadStyle = 3

$conn2_4 = $s.Shapes.AddConnector(1, 0, 0, 1, 1)
$conn2_4.Left = 394.9346850393701
$conn2_4.Top = 291.2810629921260
$conn2_4.Width = 0.0000393700787
$conn2_4.Height = 96.0000393700787
$conn2_4.Name = "Straight Arrow Connector 44"
$conn2_4.Line.Weight = 3
$conn2_4.Line.ForeColor.RGB = RGB(255,0,0)
$conn2_4.Line.EndArrowheadStyle = 3

$conn3_4 = $s.Shapes.AddLine(0, 0, 1, 1)
$conn3_4.Left = 322.9908267716535
$conn3_4.Top = 222.0000393700788
$conn3_4.Width = 0.0001181102362
$conn3_4.Height = 69.2810629921260
$conn3_4.Flip(0)
$conn3_4.Name = "Straight Connector 45"
$conn3_4.Line.ForeColor.ObjectThemeColor = 1
$conn3_4.Line.Weight = 1
$conn3_4.Line.DashStyle = 4
$conn3_4.ConnectorFormat.BeginConnect($oval4, 0)

$conn4_4 = $s.Shapes.AddLine(0, 0, 1, 1)
$conn4_4.Left = 329.9727165354331
$conn4_4.Top = 291.7107480314961
$conn4_4.Width = 67.9358661417323
$conn4_4.Height = 0.0000393700787
$conn4_4.Name = "Straight Connector 46"
$conn4_4.Line.ForeColor.ObjectThemeColor = 1
$conn4_4.Line.Weight = 1
$conn4_4.Line.DashStyle = 4

$grp4 = $s.Shapes.Range(@($rect4.Name, $oval4.Name, $conn1_4.Name, $conn2_4.Name, $conn3_4.Name, $conn4_4.Name)).Group()
$grp4.Name = "Group 40"
$grp4.Left = 410.9416141732283
$grp4.Top = 336.3942125984252
$grp4.Width = 120.0000393700787
$grp4.Height = 120.7883858267717
$grp4.Rotation = 270

